$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("indice")

# 1. Rename the long description string (held by A7) to its short variable
#    code "RTPG_PextQ" (matches the naming convention used elsewhere, e.g.
#    RTPG_mecanismos / RTPG_PextQ_vs_pib).
$ws.Range("A7").Value = "RTPG_PextQ"

# 2. Swap the contents of row 6 and row 7 (both columns A and B).
$a6 = $ws.Range("A6").Value2
$b6 = $ws.Range("B6").Value2
$a7 = $ws.Range("A7").Value2
$b7 = $ws.Range("B7").Value2

$ws.Range("A6").Value = $a7
$ws.Range("B6").Value = $b7
$ws.Range("A7").Value = $a6
$ws.Range("B7").Value = $b6

# 3. Move the selection/active cell to A7.
$ws.Range("A7").Select()

# 4. Resize the workbook window.
$excel.ActiveWindow.Width = 19800
$excel.ActiveWindow.Height = 11760
